$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. The source file stores every
# Price/Volume cell as literal text (not a native Excel number), which is
# what preserves exact formatting like trailing zeros ("305.80"), percent
# signs ("5.76%") and thousands separators ("2,458.02%"). So each cell is
# switched to Text format before the write (otherwise Excel would parse the
# numeric-looking text and silently normalize/round it), then the style is
# reset back to Normal/General afterwards to match the workbook's original,
# unstyled data cells.
$updates = [ordered]@{
    "D2" = "305.80"
    "E2" = "5.76%"
    "D3" = "32.24"
    "E3" = "9.69%"
    "D4" = "5.330"
    "E4" = "4.17%"
    "D5" = "0.07417"
    "E5" = "11.04%"
    "D6" = "7.761"
    "E6" = "5.28%"
    "D7" = "3.694"
    "E7" = "8.32%"
    "D8" = "1.538"
    "E8" = "13.37%"
    "D9" = "0.9108"
    "E9" = "-0.89%"
    "D10" = "0.01649"
    "E10" = "2,458.02%"
    "D11" = "0.1665"
    "E11" = "4.94%"
    "D12" = "0.07638"
    "E12" = "13.85%"
    "D13" = "0.07948"
    "E13" = "3.10%"
    "D14" = "0.03024"
    "E14" = "2.51%"
    "D15" = "0.09856"
    "E15" = "9.58%"
    "D16" = "0.001525"
    "E16" = "-2.77%"
    "E17" = "1.09%"
    "D18" = "0.006364"
    "E18" = "1.49%"
    "D19" = "3.494"
    "E19" = "1.37%"
    "E20" = "0.89%"
    "E21" = "1.63%"
    "D22" = "0.1336"
    "E22" = "2.04%"
    "D23" = "4.222"
    "E23" = "3.07%"
    "D24" = "0.1629"
    "E24" = "3.81%"
    "D26" = "0.004504"
    "E26" = "9.13%"
    "E27" = "-6.45%"
    "D28" = "0.0001741"
    "E28" = "7.54%"
    "D40" = "0.04491"
    "E40" = "6.56%"
    "D41" = "0.007378"
    "E41" = "9.70%"
    "D42" = "0.1360"
    "E42" = "9.61%"
    "E43" = "14.18%"
    "D44" = "0.01363"
    "E44" = "12.01%"
    "D45" = "0.00006141"
    "E45" = "7.26%"
    "D46" = "1.892"
    "E46" = "-4.02%"
    "D47" = "0.01300"
    "E47" = "-0.59%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
